$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell C5: remove the trailing "()" from "Up to 20 Stocks()"
$ws.Range("C5").Value = "Up to 20 Stocks"

# Select the edited cell, matching the saved selection state
$ws.Range("C5").Select()
